# Rename the "expert" sheets as part of breaking them out into CSVs for
# GitHub version control: fix the "sedementation" typo and make both
# sheet names lower-case to match the new csv-based workflow.
$wb = $excel.ActiveWorkbook

$sedSheet = $wb.Worksheets.Item("sedementation")
$sedSheet.Name = "sedimentation"

$eutroSheet = $wb.Worksheets.Item("Eutrophication")
$eutroSheet.Name = "eutrophication"

# Make the sedimentation sheet the active/selected tab (previously the
# eutrophication sheet was active).
$sedSheet.Activate()
